$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  45933, 4175.28711837727, 5223.43110416981, 6812, 4971.757993, -33.0040842169775),
    @(3,  45934, 824.526936987101, 2972.43031319701, 3620, 1386.369046, -3.57198240792051),
    @(4,  45935, 703.796628874127, 2560.95082234209, 3620, 1243.135093, -21.6546130638349),
    @(5,  45936, 4173.32047864009, 5215.72955659538, 3620, 4967.319353, 99.5720179564707),
    @(6,  45937, 4164.47081871551, 5013.64536330373, 3620, 4956.215963, 91.0579378161758),
    @(7,  45938, 4355.11221261204, 4756.81719354, 3620, 6266.585594, 127.012107288665),
    @(8,  45939, 4355.11221261204, 4506.80700799758, 3620, 6266.585594, 116.595016224398),
    @(9,  45940, 5445.69835475284, 4415.5175600099, 3620, 8203.301892, 148.046712385711),
    @(10, 45941, 1503.60506415849, 2402.43600475053, 3620, 3877.880684, 48.1963176913351),
    @(11, 45942, 1399.99543622148, 2346.92646387509, 3620, 3764.615186, 45.4810922355672),
    @(12, 45943, 5841.28374289435, 5212.00646101023, 3620, 8951.246298, 195.915375671495),
    @(13, 45944, 5841.28374289435, 5965.02514258658, 3620, 8951.246298, 227.29115407051),
    @(14, 45945, 5841.28374289435, 6698.75021995506, 3620, 8970.797255, 258.677655502529),
    @(15, 45946, 5841.28374289435, 6728.83195680652, 3620, 8970.797255, 259.931061204674)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
